$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D (Aquaculture production (metric tons)) entirely,
# shifting Capture fisheries production (metric tons) left from E to D.
$ws.Columns("D").Delete()

# Update selection to match the saved view state.
$ws.Range("I10").Select() | Out-Null
